$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.3
$ws.Range("V2").Value = 1.25
$ws.Range("F3").Value = 2.32
$ws.Range("G3").Value = 2.66
$ws.Range("H3").Value = 3.1
$ws.Range("L3").Value = 1.41
$ws.Range("P3").Value = 1.85
$ws.Range("Q3").Value = 1.94
$ws.Range("R3").Value = 1.32
$ws.Range("S3").Value = 3.45
$ws.Range("T3").Value = 1.64
$ws.Range("V3").Value = 1.37
$ws.Range("W3").Value = 1.6
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 15
$ws.Range("AC3").Value = 9.4
$ws.Range("AH3").Value = 22
$ws.Range("AN3").Value = 27
$ws.Range("T4").Value = 1.64
$ws.Range("U4").Value = 2.06
$ws.Range("K5").Value = 7.6
$ws.Range("U5").Value = 1.83
$ws.Range("G6").Value = 2.56
$ws.Range("I6").Value = 3.5
$ws.Range("J6").Value = 3.2
$ws.Range("W6").Value = 1.64
$ws.Range("AN6").Value = 27
$ws.Range("F7").Value = 1.49
$ws.Range("I7").Value = 8.4
$ws.Range("K7").Value = 6
$ws.Range("P7").Value = 2.58
$ws.Range("V7").Value = 1.16
$ws.Range("W7").Value = 2.66
$ws.Range("Z7").Value = 70
$ws.Range("AD7").Value = 27
$ws.Range("AH7").Value = 24
$ws.Range("AI7").Value = 80
$ws.Range("AO7").Value = 85
$ws.Range("S8").Value = 2.34
$ws.Range("H9").Value = 1.8
$ws.Range("I9").Value = 1.92
$ws.Range("J9").Value = 3.5
$ws.Range("S9").Value = 3.7
$ws.Range("F10").Value = 7.6
$ws.Range("G10").Value = 12
$ws.Range("H10").Value = 1.32
$ws.Range("I10").Value = 1.42
$ws.Range("J10").Value = 5.4
$ws.Range("K10").Value = 7.4
$ws.Range("N10").Value = 5.5
$ws.Range("O10").Value = 1.17
$ws.Range("Q10").Value = 1.47
$ws.Range("R10").Value = 1.63
$ws.Range("V10").Value = 3.35
$ws.Range("AI10").Value = 32
$ws.Range("F11").Value = 1.53
$ws.Range("I11").Value = 9.4
$ws.Range("K11").Value = 4.4
$ws.Range("N11").Value = 2.98
$ws.Range("P11").Value = 1.69
$ws.Range("T11").Value = 2.4
$ws.Range("U11").Value = 1.62
$ws.Range("AH11").Value = 36
$ws.Range("AI11").Value = 210
$ws.Range("F14").Value = 1.44
$ws.Range("G14").Value = 1.53
$ws.Range("J14").Value = 3.85
$ws.Range("K14").Value = 5.3
$ws.Range("F16").Value = 1.56
$ws.Range("G16").Value = 1.71
$ws.Range("J16").Value = 3.75
$ws.Range("N16").Value = 3.55
$ws.Range("O16").Value = 1.29
$ws.Range("P16").Value = 1.86
$ws.Range("Q16").Value = 1.89
$ws.Range("S17").Value = 4.6
$ws.Range("T17").Value = 2
$ws.Range("T18").Value = 1.8
$ws.Range("P19").Value = 1.71
$ws.Range("Q19").Value = 2.38
$ws.Range("R19").Value = 1.26
$ws.Range("AL19").Value = 50
$ws.Range("F20").Value = 1.93
$ws.Range("G20").Value = 1.99
$ws.Range("J20").Value = 3.7
$ws.Range("K20").Value = 3.9
$ws.Range("P20").Value = 2.04
$ws.Range("Q20").Value = 1.79
$ws.Range("V20").Value = 1.29
$ws.Range("W20").Value = 2
$ws.Range("F21").Value = 1.86
$ws.Range("G21").Value = 1.87
$ws.Range("P21").Value = 2.6
$ws.Range("S21").Value = 2.48
$ws.Range("U21").Value = 2.58
$ws.Range("W21").Value = 2.14
$ws.Range("AG21").Value = 9.800000000000001
$ws.Range("AJ21").Value = 21
$ws.Range("AK21").Value = 16
$ws.Range("AL21").Value = 24
$ws.Range("AN21").Value = 8.199999999999999
$ws.Range("AO21").Value = 30
$ws.Range("F22").Value = 2.68
$ws.Range("M22").Value = 1.07
$ws.Range("R22").Value = 1.35
$ws.Range("S22").Value = 3.65
$ws.Range("Y22").Value = 11.5
$ws.Range("AL22").Value = 44
$ws.Range("F23").Value = 1.2
$ws.Range("G23").Value = 1.21
$ws.Range("J23").Value = 7.8
$ws.Range("K23").Value = 8.4
$ws.Range("N23").Value = 5.1
$ws.Range("P23").Value = 2.4
$ws.Range("Q23").Value = 1.63
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.6
$ws.Range("U23").Value = 1.59
$ws.Range("W23").Value = 5.7
$ws.Range("Z23").Value = 260
$ws.Range("AE23").Value = 550
$ws.Range("AI23").Value = 370
$ws.Range("F24").Value = 3.3
$ws.Range("G24").Value = 4.1
$ws.Range("H24").Value = 1.93
$ws.Range("I24").Value = 2.08
$ws.Range("K24").Value = 4.7
$ws.Range("P24").Value = 2.22
$ws.Range("T24").Value = 1.63
$ws.Range("U24").Value = 2.22
$ws.Range("V24").Value = 1.92
$ws.Range("W24").Value = 1.33
$ws.Range("AG24").Value = 16
$ws.Range("AK24").Value = 44
$ws.Range("AM24").Value = 80
